$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 90 (pushes existing rows 90+ down to 91+)
$ws.Rows.Item(90).Insert()

# Resize Table1 to include the new row (table now spans A8:K151)
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K151"))

# Copy cell formatting (borders/number format) from the existing "year header" row (69)
# into the new row 90, restricted to the table's columns (A:K) only
$ws.Range("A69:K69").Copy()
$ws.Range("A90:K90").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 90 is the "2024" year header
$ws.Range("A90").Value() = "2024"

# Row 87: SL(1-0-0) leave record gets 1.25 days earned
$ws.Range("C87").Value() = 1.25

# Row 89: new leave entry "SP(2-0-0)" for 1.25 days, remark period 01/02,04/2024
$ws.Range("B89").Value() = "SP(2-0-0)"
$ws.Range("C89").Value() = 1.25
$ws.Range("K89").Value() = "01/02,04/2024"

# Row 91 (previously row 90 before insert): SL(1-0-0) with 1 day absence undertime w/o pay, remark date
$ws.Range("B91").Value() = "SL(1-0-0)"
$ws.Range("H91").Value() = 1
$ws.Range("K91").Value() = 45289

Write-Host "done"
